$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab/title to reflect the new "through" date
$ws.Name = "Through 2022-07-15"

# Update the label in column A for the partial-July row (row 8)
$ws.Range("A8").Value = "July (through 07-15)"

# Update June 2022 figure (row 7 = June row), which shifted slightly with the new data pull
$ws.Range("I7").Value = 142

# Update July row (row 8) values across years
$ws.Range("B8").Value = 21
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 29
$ws.Range("E8").Value = 36
$ws.Range("G8").Value = 57
$ws.Range("H8").Value = 73
$ws.Range("I8").Value = 82

# Update Total row (row 9) values across years
$ws.Range("B9").Value = 146
$ws.Range("C9").Value = 278
$ws.Range("D9").Value = 419
$ws.Range("E9").Value = 389
$ws.Range("G9").Value = 529
$ws.Range("H9").Value = 833
$ws.Range("I9").Value = 887
